# First wave of updates from institute PIs
# - Add William Shields' ORCID (row 3, column N "ORCID ") which was
#   previously blank.
# - Leave the cursor/selection on the newly-edited cell, as the author did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N3").Value = "0000-0002-4825-0303"

# Move the active selection to the cell that was just edited.
$ws.Range("N3").Select()
